$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2 through 13
# from 46074 (2026-02-21) to 46075 (2026-02-22)
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 46075
}
